$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 269.44446
$ws.Range("J28").Value = 286.25
$ws.Range("L28").Value = 286.25
$ws.Range("N28").Value = -1256.25
$ws.Range("H70").Value = 3011.8
$ws.Range("I70").Value = 2748
$ws.Range("J70").Value = 3187.6667
$ws.Range("K70").Value = 8244
$ws.Range("L70").Value = 9563.000100000001
$ws.Range("M70").Value = -7974
$ws.Range("N70").Value = -10103.0001
$ws.Range("H73").Value = 3011.8
$ws.Range("I73").Value = 2748
$ws.Range("J73").Value = 3187.6667
$ws.Range("K73").Value = 8244
$ws.Range("L73").Value = 9563.000100000001
$ws.Range("M73").Value = -7308
$ws.Range("N73").Value = -11435.0001
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H98").Value = 3945.85
$ws.Range("I98").Value = 4312.788
$ws.Range("K98").Value = 4312.788
$ws.Range("M98").Value = -2814.788
$ws.Range("H122").Value = 3945.85
$ws.Range("I122").Value = 4312.788
$ws.Range("K122").Value = 12938.364
$ws.Range("M122").Value = -10488.364
$ws.Range("H132").Value = 3248977.8
$ws.Range("I132").Value = 3573445.2
$ws.Range("K132").Value = 10720335.6
$ws.Range("M132").Value = -10717805.6
$ws.Range("H137").Value = 1415.4814
$ws.Range("J137").Value = 1624.875
$ws.Range("L137").Value = 4874.625
$ws.Range("N137").Value = -9974.625
$ws.Range("H138").Value = 1660.7407
$ws.Range("I138").Value = 928.3871
$ws.Range("J138").Value = 2647.8262
$ws.Range("K138").Value = 2785.1613
$ws.Range("L138").Value = 7943.4786
$ws.Range("M138").Value = 2354.8387
$ws.Range("N138").Value = -18223.4786

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18700.016
$ws.Range("I32").Value = 20027.21
$ws.Range("K32").Value = 20027.21
$ws.Range("M32").Value = -19740.21

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 74943.71000000001
$ws.Range("I134").Value = 115023.555
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 345070.665
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -342535.665
$ws.Range("N134").Value = -13470

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1654463.8
$ws.Range("I31").Value = 1930.0333
$ws.Range("K31").Value = 1930.0333
$ws.Range("M31").Value = -1635.0333
$ws.Range("H34").Value = 1654463.8
$ws.Range("I34").Value = 1930.0333
$ws.Range("K34").Value = 1930.0333
$ws.Range("M34").Value = -1728.0333
$ws.Range("H58").Value = 1333.9445
$ws.Range("I58").Value = 1301
$ws.Range("J58").Value = 1385.7142
$ws.Range("K58").Value = 1301
$ws.Range("L58").Value = 1385.7142
$ws.Range("M58").Value = -1098
$ws.Range("N58").Value = -1791.7142
$ws.Range("H99").Value = 2412.7
$ws.Range("J99").Value = 2447.4443
$ws.Range("L99").Value = 2447.4443
$ws.Range("N99").Value = -5443.4443
$ws.Range("H126").Value = 2412.7
$ws.Range("J126").Value = 2447.4443
$ws.Range("L126").Value = 7342.3329
$ws.Range("N126").Value = -12282.3329
$ws.Range("H136").Value = 1333.9445
$ws.Range("I136").Value = 1301
$ws.Range("J136").Value = 1385.7142
$ws.Range("K136").Value = 3903
$ws.Range("L136").Value = 4157.142599999999
$ws.Range("M136").Value = -1353
$ws.Range("N136").Value = -9257.142599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 158.44444
$ws.Range("I8").Value = 158.44444
$ws.Range("K8").Value = 475.33332
$ws.Range("M8").Value = -336.33332
$ws.Range("H68").Value = 1086.3191
$ws.Range("I68").Value = 937.09375
$ws.Range("J68").Value = 1404.6666
$ws.Range("K68").Value = 2811.28125
$ws.Range("L68").Value = 4213.9998
$ws.Range("M68").Value = -2000.28125
$ws.Range("N68").Value = -5835.9998
$ws.Range("H71").Value = 1086.3191
$ws.Range("I71").Value = 937.09375
$ws.Range("J71").Value = 1404.6666
$ws.Range("K71").Value = 8433.84375
$ws.Range("L71").Value = 12641.9994
$ws.Range("M71").Value = -4377.84375
$ws.Range("N71").Value = -20753.9994
$ws.Range("H131").Value = 29708
$ws.Range("J131").Value = 32091.375
$ws.Range("L131").Value = 96274.125
$ws.Range("N131").Value = -106354.125

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 12103.9
$ws.Range("I126").Value = 7437.8
$ws.Range("K126").Value = 22313.4
$ws.Range("M126").Value = -19843.4
$ws.Range("H132").Value = 31949
$ws.Range("I132").Value = 45370.695
$ws.Range("K132").Value = 136112.085
$ws.Range("M132").Value = -133582.085

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H40").Value = 2527.182
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2966.5
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2966.5
$ws.Range("M40").Value = -1864
$ws.Range("N40").Value = -3238.5
$ws.Range("H132").Value = 5448.378
$ws.Range("I132").Value = 7346
$ws.Range("J132").Value = 2322.8823
$ws.Range("K132").Value = 22038
$ws.Range("L132").Value = 6968.646900000001
$ws.Range("M132").Value = -19508
$ws.Range("N132").Value = -12028.6469
$ws.Range("H136").Value = 4531.4116
$ws.Range("I136").Value = 6518.278
$ws.Range("J136").Value = 2296.1875
$ws.Range("K136").Value = 19554.834
$ws.Range("L136").Value = 6888.5625
$ws.Range("M136").Value = -17004.834
$ws.Range("N136").Value = -11988.5625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4088.0908
$ws.Range("I122").Value = 5834
$ws.Range("J122").Value = 1993
$ws.Range("K122").Value = 17502
$ws.Range("L122").Value = 5979
$ws.Range("M122").Value = -15052
$ws.Range("N122").Value = -10879
$ws.Range("H132").Value = 1467.7142
$ws.Range("I132").Value = 1042.4166
$ws.Range("J132").Value = 1786.6875
$ws.Range("K132").Value = 3127.2498
$ws.Range("L132").Value = 5360.0625
$ws.Range("M132").Value = -597.2498000000001
$ws.Range("N132").Value = -10420.0625
$ws.Range("H136").Value = 1286.3549
$ws.Range("I136").Value = 1401.4546
$ws.Range("K136").Value = 4204.3638
$ws.Range("M136").Value = -1654.3638
